# Generate Report for handoff
#
# The "735e2a3e-aa66-4b58-9612-6255bd4eaca9.md" source file was re-handed-off
# under a new id ("572775a6-2d63-486d-8196-f83075dc5894.md") with a new
# transform hash ("6f0063a7efe1d520d593ca110b4191ce936548d6") and new handoff
# timestamps; the old (now-stale) "a8d7baa1-8806-428c-92bd-f19d05061f94.md"
# "Handoff transform failed" row is dropped from every sheet.

$wb = $excel.ActiveWorkbook

$oldGuid = "735e2a3e-aa66-4b58-9612-6255bd4eaca9"
$newGuid = "572775a6-2d63-486d-8196-f83075dc5894"

$oldHash = "cab1c72d2cf683739f2b8d96785b874c26c91cb3"
$newHash = "6f0063a7efe1d520d593ca110b4191ce936548d6"

$oldZhDate = "2016-02-19 06:04:09"
$newZhDate = "2016-02-19 06:04:58"

$oldDeDate = "2016-02-19 06:04:22"
$newDeDate = "2016-02-19 06:05:11"

$oldMd = $oldGuid + ".md"
$newMd = $newGuid + ".md"

$oldZhXlf = $oldGuid + "." + $oldHash + ".zh-cn.xlf"
$newZhXlf = $newGuid + "." + $newHash + ".zh-cn.xlf"

$oldDeXlf = $oldGuid + "." + $oldHash + ".de-de.xlf"
$newDeXlf = $newGuid + "." + $newHash + ".de-de.xlf"

$sheetNames = @("Overview", "zh-cn", "de-de")

# 1) Drop the stale "a8d7baa1-....md" / "Handoff transform failed" row
#    (row 3 on every sheet) - shifts the ".localization-config" row up.
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(3).Delete()
}

# 2) Rename the re-handed-off file + refresh its handoff artifacts/timestamps
#    everywhere it is referenced.
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($r = 1; $r -le 3; $r++) {
        for ($c = 1; $c -le 9; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null) {
                $s = "" + $val
                $new = $s
                if ($s -eq $oldMd) {
                    $new = $newMd
                } elseif ($s -eq $oldZhXlf) {
                    $new = $newZhXlf
                } elseif ($s -eq $oldDeXlf) {
                    $new = $newDeXlf
                } elseif ($s -eq $oldZhDate) {
                    $new = $newZhDate
                } elseif ($s -eq $oldDeDate) {
                    $new = $newDeDate
                }
                if ($new -ne $s) {
                    $cell.Value = $new
                }
            }
        }
    }
}

# 3) Rebuild hyperlinks per-sheet (row delete above does not renumber /
#    retarget the surviving hyperlinks, and the renamed file needs a new
#    target URL too).
$baseMd = "https://github.com/OpenLocalizationTest/oltest/blob/9dedf9c27530375a8e84645925f3bf59948bdc39/e2e/"
$baseCfg = "https://github.com/OpenLocalizationTest/oltest/blob/9dedf9c27530375a8e84645925f3bf59948bdc39/.localization-config"
$baseZhXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36736f186b9123ad4dd12f51cc9620b165ddd5c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/"
$baseDeXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2382c85a6bc1f1ecde69102d4ae8b5a6fefef948/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($baseMd + $newMd), "", "", $newMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $baseCfg, "", "", ".localization-config")

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($baseMd + $newMd), "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), ($baseZhXlf + $newZhXlf), "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $baseCfg, "", "", ".localization-config")

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($baseMd + $newMd), "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), ($baseDeXlf + $newDeXlf), "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $baseCfg, "", "", ".localization-config")

# 4) `Hyperlinks.Add` re-styles touched cells with Excel's generic built-in
#    "Hyperlink" look; restore the workbook's own custom hyperlink font
#    (single underline, #6495ED) on every linked cell so the restyle is a
#    no-op relative to the original workbook.
$linkColor = 15570276  # BGR packing of RGB 6495ED, as Font.Color expects
$hyperlinkCells = @(
    $wsOverview.Range("A2"), $wsOverview.Range("A3"),
    $wsZh.Range("A2"), $wsZh.Range("C2"), $wsZh.Range("A3"),
    $wsDe.Range("A2"), $wsDe.Range("C2"), $wsDe.Range("A3")
)
foreach ($cell in $hyperlinkCells) {
    $cell.Font.Underline = 2
    $cell.Font.Color = $linkColor
}
